$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 45 corresponds to MindSensorPressureSensor - add Developer/finished/Mode/interface
$ws.Range("D45").Value = "Lawrie"
$ws.Range("E45").Value = "N"
$ws.Range("F45").Value = "Pressure"
$ws.Range("G45").Value = "SampleProvider"

# Row 69 corresponds to RCXTemperatureSensor - add Developer/finished/Mode/interface
$ws.Range("D69").Value = "Lawrie"
$ws.Range("E69").Value = "N"
$ws.Range("F69").Value = "Temperature"
$ws.Range("G69").Value = "SampleProvider"

# Update view: frozen pane top-left cell and active selection
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("I66").Select()
